$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 0.6973690231068405
$ws.Range("K2").Value = 0.7208591480536559
$ws.Range("L2").Value = 0.6973690231068405
$ws.Range("M2").Value = 0.6953594983715047
$ws.Range("R2").Value = 0.7660032029283916
$ws.Range("S2").Value = 0.7760820378100716
$ws.Range("T2").Value = 0.7660032029283916
$ws.Range("U2").Value = 0.7654431671124573
$ws.Range("W2").Value = 0.7877672318690958
$ws.Range("Y2").Value = 0.7772809472569253
$ws.Range("F3").Value = 0.8133379089453214
$ws.Range("G3").Value = 0.8218721920011042
$ws.Range("H3").Value = 0.8133379089453214
$ws.Range("I3").Value = 0.8147808147714887
$ws.Range("K3").Value = 0.8331590494638522
$ws.Range("M3").Value = 0.8247243502739942
$ws.Range("R3").Value = 0.8519789521848548
$ws.Range("S3").Value = 0.8598142707282346
$ws.Range("T3").Value = 0.8519789521848548
$ws.Range("U3").Value = 0.852081097935892
$ws.Range("V3").Value = 0.8498284145504462
$ws.Range("W3").Value = 0.8581762562608208
$ws.Range("X3").Value = 0.8498284145504462
$ws.Range("Y3").Value = 0.8496333227015507
$ws.Range("F4").Value = 0.8734385724090596
$ws.Range("G4").Value = 0.8801279785676346
$ws.Range("H4").Value = 0.8734385724090596
$ws.Range("I4").Value = 0.8738176868073346
$ws.Range("J4").Value = 0.8433996797071609
$ws.Range("K4").Value = 0.8505098024431348
$ws.Range("L4").Value = 0.8433996797071609
$ws.Range("M4").Value = 0.8440484512013443
$ws.Range("R4").Value = 0.8606039807824297
$ws.Range("S4").Value = 0.8652903416814567
$ws.Range("T4").Value = 0.8606039807824297
$ws.Range("U4").Value = 0.8600211334618383
$ws.Range("V4").Value = 0.8605811027224892
$ws.Range("W4").Value = 0.8660553055821871
$ws.Range("X4").Value = 0.8605811027224892
$ws.Range("Y4").Value = 0.8599587702034434
$ws.Range("N5").Value = 0.8369251887439946
$ws.Range("O5").Value = 0.8474318443561014
$ws.Range("P5").Value = 0.8369251887439946
$ws.Range("Q5").Value = 0.8351730171881542
$ws.Range("C6").Value = 0.8595629681419273
$ws.Range("E6").Value = 0.8539852933423724
$ws.Range("F6").Value = 0.8713109128345916
$ws.Range("G6").Value = 0.8755292665935478
$ws.Range("H6").Value = 0.8713109128345916
$ws.Range("I6").Value = 0.8709376453541209
$ws.Range("J6").Value = 0.8304735758407688
$ws.Range("K6").Value = 0.8403918233504697
$ws.Range("L6").Value = 0.8304735758407688
$ws.Range("M6").Value = 0.8307894683990016
$ws.Range("N6").Value = 0.8519331960649736
$ws.Range("O6").Value = 0.8577064876295909
$ws.Range("P6").Value = 0.8519331960649736
$ws.Range("Q6").Value = 0.851240140726601
$ws.Range("R6").Value = 0.8691603752001831
$ws.Range("S6").Value = 0.874993950629074
$ws.Range("T6").Value = 0.8691603752001831
$ws.Range("U6").Value = 0.8685697322098177
$ws.Range("V6").Value = 0.8540379775795014
$ws.Range("W6").Value = 0.8589613235760376
$ws.Range("X6").Value = 0.8540379775795014
$ws.Range("Y6").Value = 0.853904858120757
